$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 2 (was old row 3: FAPs / Agt / Lrp2 / MuSCs) with updated TPM-derived values
$ws.Range("A2").Value = "FAPs"
$ws.Range("B2").Value = "Agt"
$ws.Range("C2").Value = "Lrp2"
$ws.Range("D2").Value = "MuSCs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.8072116666666668
$ws.Range("H2").Value = 2.421635
$ws.Range("I2").Value = 0.7485686721305167
$ws.Range("J2").Value = 0.8170456278433896
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.032708
$ws.Range("N2").Value = 0.065416
$ws.Range("O2").Value = 1
$ws.Range("P2").Value = 1
$ws.Range("Q2").Value = 0.02640227919333334
$ws.Range("R2").Value = 0.15841367516
$ws.Range("S2").Value = 0.7485686721305167
$ws.Range("T2").Value = 0.8170456278433896

# New row 3 (was old row 4: MuSCs / Agt / Lrp2 / MuSCs) with updated TPM-derived values
$ws.Range("A3").Value = "MuSCs"
$ws.Range("B3").Value = "Agt"
$ws.Range("C3").Value = "Lrp2"
$ws.Range("D3").Value = "MuSCs"
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0.2711285
$ws.Range("H3").Value = 0.542257
$ws.Range("I3").Value = 0.2514313278694834
$ws.Range("J3").Value = 0.1829543721566103
$ws.Range("K3").Value = 2
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 0.032708
$ws.Range("N3").Value = 0.065416
$ws.Range("O3").Value = 1
$ws.Range("P3").Value = 1
$ws.Range("Q3").Value = 0.008868070978000001
$ws.Range("R3").Value = 0.035472283912
$ws.Range("S3").Value = 0.2514313278694834
$ws.Range("T3").Value = 0.1829543721566103

# Remove old rows 4 and 5 (ECs/Neutrophils data no longer present)
$ws.Range("A4:T5").EntireRow.Delete()
